# Repeat-trader run: fill in today's PriceChange/UpDown for the prior row
# and append the newest row of scored data (row 10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- complete row 9 (PriceChange / UpDown were computed a day later) ---
$ws.Range("X9").Value = -1.2099989999999963
$ws.Range("Y9").Value = "Down"

# --- append row 10 ---
$ws.Range("A10").Value = 42653.879016203704
$ws.Range("B10").Value = 2
$ws.Range("C10").Value = "Neutral"
$ws.Range("D10").Value = 8
$ws.Range("E10").Value = 17061
$ws.Range("F10").Value = 2705
$ws.Range("G10").Value = 56
$ws.Range("H10").Value = 40
$ws.Range("I10").Value = 86
$ws.Range("J10").Value = 12
$ws.Range("K10").Value = 22291
$ws.Range("L10").Value = 350
$ws.Range("M10").Value = 250
$ws.Range("N10").Value = 122
$ws.Range("O10").Value = 18
$ws.Range("P10").Value = "Noun"
$ws.Range("Q10").Value = 35.550971360736582
$ws.Range("R10").Value = -24.44

$ws.Range("S10").NumberFormat = "0.00%"
$ws.Range("S10").Value = -0.1153
$ws.Range("T10").NumberFormat = "0.00%"
$ws.Range("T10").Value = -0.047

$ws.Range("U10").Value = 6.45
$ws.Range("V10").Value = 1.88
$ws.Range("W10").Value = 1

# --- re-fit the bestFit columns now that the data (and its widest values) changed ---
$ws.Columns.Item(1).ColumnWidth = 14.5
$ws.Columns.Item(2).ColumnWidth = 7.666666666666667
$ws.Columns.Item(3).ColumnWidth = 5.666666666666667
$ws.Columns.Item(4).ColumnWidth = 11.333333333333334
$ws.Columns.Item(5).ColumnWidth = 8.666666666666666
$ws.Columns.Item(6).ColumnWidth = 11.5
$ws.Columns.Item(7).ColumnWidth = 18.5
$ws.Columns.Item(8).ColumnWidth = 18.5
$ws.Columns.Item(9).ColumnWidth = 19.5
$ws.Columns.Item(10).ColumnWidth = 19.833333333333332
$ws.Columns.Item(11).ColumnWidth = 9.5
$ws.Columns.Item(12).ColumnWidth = 13.5
$ws.Columns.Item(13).ColumnWidth = 13.833333333333334
